# edit.ps1
# Scheduled-runner data refresh for Sheets/Midgardsormr_Profits.xlsx
# Updates cached market-board / leve-profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 3
$ws.Range("H3").Value = 49998
$ws.Range("J3").Value = 49998
$ws.Range("L3").Value = 49998
$ws.Range("N3").Value = -50226
# row 15
$ws.Range("H15").Value = 614.5952
$ws.Range("I15").Value = 614.5952
$ws.Range("K15").Value = 1843.7856
$ws.Range("M15").Value = -1674.7856
# row 102
$ws.Range("H102").Value = 49998
$ws.Range("J102").Value = 49998
$ws.Range("L102").Value = 49998
$ws.Range("N102").Value = -56488
# row 107
$ws.Range("H107").Value = 1566.5625
$ws.Range("I107").Value = 1566.5625
$ws.Range("K107").Value = 1566.5625
$ws.Range("M107").Value = 353.4375
# row 111
$ws.Range("H111").Value = 1543.3125
$ws.Range("I111").Value = 1459.875
$ws.Range("J111").Value = 1626.75
$ws.Range("K111").Value = 4379.625
$ws.Range("L111").Value = 4880.25
$ws.Range("M111").Value = -1312.625
$ws.Range("N111").Value = -11014.25
# row 112
$ws.Range("H112").Value = 5897.7036
$ws.Range("I112").Value = 899.5
$ws.Range("J112").Value = 6766.9565
$ws.Range("K112").Value = 2698.5
$ws.Range("L112").Value = 20300.8695
$ws.Range("M112").Value = -1590.5
$ws.Range("N112").Value = -22516.8695
# row 113
$ws.Range("H113").Value = 6111
$ws.Range("J113").Value = 6333.3335
$ws.Range("L113").Value = 6333.3335
$ws.Range("N113").Value = -12841.3335
# row 132
$ws.Range("H132").Value = 2690.2778
$ws.Range("I132").Value = 2173.5
$ws.Range("J132").Value = 3723.8333
$ws.Range("K132").Value = 6520.5
$ws.Range("L132").Value = 11171.4999
$ws.Range("M132").Value = -3990.5
$ws.Range("N132").Value = -16231.4999
# row 138
$ws.Range("H138").Value = 1530576.8
$ws.Range("I138").Value = 1846.6875
$ws.Range("J138").Value = 2229424.8
$ws.Range("K138").Value = 5540.0625
$ws.Range("L138").Value = 6688274.399999999
$ws.Range("M138").Value = -400.0625
$ws.Range("N138").Value = -6698554.399999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("H61").Value = 6784.75
$ws.Range("I61").Value = 4191.143
$ws.Range("J61").Value = 15862.375
$ws.Range("K61").Value = 4191.143
$ws.Range("L61").Value = 15862.375
$ws.Range("M61").Value = -3979.143
$ws.Range("N61").Value = -16286.375
# row 74
$ws.Range("H74").Value = 3396.4
$ws.Range("I74").Value = 1194.8857
$ws.Range("J74").Value = 8533.267
$ws.Range("K74").Value = 1194.8857
$ws.Range("L74").Value = 8533.267
$ws.Range("M74").Value = -320.8857
$ws.Range("N74").Value = -10281.267
# row 77
$ws.Range("H77").Value = 3396.4
$ws.Range("I77").Value = 1194.8857
$ws.Range("J77").Value = 8533.267
$ws.Range("K77").Value = 5974.4285
$ws.Range("L77").Value = 42666.335
$ws.Range("M77").Value = -1606.4285
$ws.Range("N77").Value = -51402.335
# row 102
$ws.Range("H102").Value = 1560.7826
$ws.Range("I102").Value = 1502.5
$ws.Range("K102").Value = 1502.5
$ws.Range("M102").Value = 119.5
# row 132
$ws.Range("H132").Value = 4877.773
$ws.Range("I132").Value = 4851.875
$ws.Range("J132").Value = 4946.8335
$ws.Range("K132").Value = 14555.625
$ws.Range("L132").Value = 14840.5005
$ws.Range("M132").Value = -12025.625
$ws.Range("N132").Value = -19900.5005
# row 135
$ws.Range("H135").Value = 87829.8
$ws.Range("J135").Value = 87829.8
$ws.Range("L135").Value = 87829.8
$ws.Range("N135").Value = -97969.8
# row 136
$ws.Range("H136").Value = 6784.75
$ws.Range("I136").Value = 4191.143
$ws.Range("J136").Value = 15862.375
$ws.Range("K136").Value = 12573.429
$ws.Range("L136").Value = 47587.125
$ws.Range("M136").Value = -10023.429
$ws.Range("N136").Value = -52687.125
# row 139
$ws.Range("H139").Value = 104420
$ws.Range("J139").Value = 104420
$ws.Range("L139").Value = 104420
$ws.Range("N139").Value = -114700

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 69249.836
$ws.Range("I20").Value = 63099.8
$ws.Range("K20").Value = 63099.8
$ws.Range("M20").Value = -62852.8
# row 105
$ws.Range("H105").Value = 4977.4546
$ws.Range("I105").Value = 6559
$ws.Range("K105").Value = 6559
$ws.Range("M105").Value = -4812
# row 134
$ws.Range("H134").Value = 4022.8215
$ws.Range("I134").Value = 3880.9167
$ws.Range("K134").Value = 11642.7501
$ws.Range("M134").Value = -9107.750100000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 19
$ws.Range("H19").Value = 791.9
$ws.Range("I19").Value = 813.6667
$ws.Range("J19").Value = 596
$ws.Range("K19").Value = 813.6667
$ws.Range("L19").Value = 596
$ws.Range("M19").Value = -643.6667
$ws.Range("N19").Value = -936
# row 24
$ws.Range("H24").Value = 791.9
$ws.Range("I24").Value = 813.6667
$ws.Range("J24").Value = 596
$ws.Range("K24").Value = 813.6667
$ws.Range("L24").Value = 596
$ws.Range("M24").Value = -643.6667
$ws.Range("N24").Value = -936
# row 31
$ws.Range("H31").Value = 2780663
$ws.Range("I31").Value = 5558327
$ws.Range("J31").Value = 2999.2778
$ws.Range("K31").Value = 5558327
$ws.Range("L31").Value = 2999.2778
$ws.Range("M31").Value = -5558032
$ws.Range("N31").Value = -3589.2778
# row 34
$ws.Range("H34").Value = 2780663
$ws.Range("I34").Value = 5558327
$ws.Range("J34").Value = 2999.2778
$ws.Range("K34").Value = 5558327
$ws.Range("L34").Value = 2999.2778
$ws.Range("M34").Value = -5558125
$ws.Range("N34").Value = -3403.2778
# row 99
$ws.Range("H99").Value = 4527.8667
$ws.Range("J99").Value = 4547.5454
$ws.Range("L99").Value = 4547.5454
$ws.Range("N99").Value = -7543.5454
# row 104
$ws.Range("H104").Value = 19500
$ws.Range("J104").Value = 19500
$ws.Range("L104").Value = 19500
$ws.Range("N104").Value = -24742
# row 105
$ws.Range("H105").Value = 4374.75
$ws.Range("I105").Value = 3999.6
$ws.Range("K105").Value = 3999.6
$ws.Range("M105").Value = -2252.6
# row 126
$ws.Range("H126").Value = 4527.8667
$ws.Range("J126").Value = 4547.5454
$ws.Range("L126").Value = 13642.6362
$ws.Range("N126").Value = -18582.6362
# row 132
$ws.Range("H132").Value = 3928.6
$ws.Range("J132").Value = 4409.25
$ws.Range("L132").Value = 13227.75
$ws.Range("N132").Value = -18287.75
# row 134
$ws.Range("H134").Value = 2705.5881
$ws.Range("I134").Value = 2080.6667
$ws.Range("J134").Value = 4205.4
$ws.Range("K134").Value = 6242.000100000001
$ws.Range("L134").Value = 12616.2
$ws.Range("M134").Value = -3707.000100000001
$ws.Range("N134").Value = -17686.2

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 752.1429000000001
$ws.Range("I5").Value = 726.8182
$ws.Range("K5").Value = 2180.4546
$ws.Range("M5").Value = -2068.4546
# row 112
$ws.Range("H112").Value = 5333.3335
$ws.Range("I112").Value = 3000
$ws.Range("K112").Value = 9000
$ws.Range("M112").Value = -7892
# row 116
$ws.Range("H116").Value = 9503.093999999999
$ws.Range("J116").Value = 9774.178
$ws.Range("L116").Value = 29322.534
$ws.Range("N116").Value = -36206.534
# row 135
$ws.Range("H135").Value = 752.1429000000001
$ws.Range("I135").Value = 726.8182
$ws.Range("K135").Value = 6541.3638
$ws.Range("M135").Value = -4006.3638

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 25
$ws.Range("H25").Value = 9909
$ws.Range("J25").Value = 9909
$ws.Range("L25").Value = 9909
$ws.Range("N25").Value = -10967
# row 80
$ws.Range("H80").Value = 7199.5
$ws.Range("I80").Value = 9199.200000000001
$ws.Range("K80").Value = 9199.200000000001
$ws.Range("M80").Value = -8201.200000000001
# row 83
$ws.Range("H83").Value = 7199.5
$ws.Range("I83").Value = 9199.200000000001
$ws.Range("K83").Value = 45996
$ws.Range("M83").Value = -41004

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 16
$ws.Range("H16").Value = 1059.3077
$ws.Range("I16").Value = 981.9
$ws.Range("J16").Value = 1317.3334
$ws.Range("K16").Value = 981.9
$ws.Range("L16").Value = 1317.3334
$ws.Range("M16").Value = -811.9
$ws.Range("N16").Value = -1657.3334
# row 30
$ws.Range("H30").Value = 4356.2
$ws.Range("I30").Value = 4356.2
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 4356.2
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -4248.2
$ws.Range("N30").ClearContents()
# row 61
$ws.Range("H61").Value = 1442
$ws.Range("I61").Value = 1789.375
$ws.Range("K61").Value = 1789.375
$ws.Range("M61").Value = -1587.375
# row 113
$ws.Range("H113").Value = 1442
$ws.Range("I113").Value = 1789.375
$ws.Range("K113").Value = 1789.375
$ws.Range("M113").Value = 380.625
# row 136
$ws.Range("H136").Value = 5088.7617
$ws.Range("I136").Value = 3817.125
$ws.Range("J136").Value = 5871.3076
$ws.Range("K136").Value = 11451.375
$ws.Range("L136").Value = 17613.9228
$ws.Range("M136").Value = -8901.375
$ws.Range("N136").Value = -22713.9228

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 12267.462
$ws.Range("I62").Value = 12082
$ws.Range("K62").Value = 12082
$ws.Range("M62").Value = -11458
# row 65
$ws.Range("H65").Value = 12267.462
$ws.Range("I65").Value = 12082
$ws.Range("K65").Value = 60410
$ws.Range("M65").Value = -57290
# row 102
$ws.Range("H102").Value = 67500
$ws.Range("J102").Value = 67500
$ws.Range("L102").Value = 67500
$ws.Range("N102").Value = -73990

